# Wolf-Fuels.xlsx update:
#  - Kudzai Tivatye's email (B2) changed from nicolaskudzai696@gmail.com
#    to kudziet221@gmail.com (old shared-string entry removed, new one
#    added at the end of the shared-strings table).
#  - B2 picks up its own distinct cell format record (still "no fill",
#    visually unchanged) instead of sharing the default style.
#  - Active selection moved from F5 to C13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email address for Kudzai Tivatye.
$ws.Range("B2").Value = "kudziet221@gmail.com"

# Touch the Interior formatting on B2 so it gets its own cell format
# record (fill stays "no fill", matching the surrounding cells visually).
$ws.Range("B2").Interior.Color = 255
$ws.Range("B2").Interior.Pattern = -4142

# Move the selected/active cell to C13.
[void]$ws.Range("C13").Select()
